$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every timestamp in column A (rows 2 through 97) forward by exactly
# one day, preserving the time-of-day fraction.
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 1
}

# Update the Actual Production values for the newly shifted day's
# morning ramp-up rows (model retrained -> lower predicted output).
$ws.Cells.Item(28, 2).Value = 7
$ws.Cells.Item(29, 2).Value = 21
$ws.Cells.Item(30, 2).Value = 51
$ws.Cells.Item(31, 2).Value = 82
